$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=45694.49027777778; B=26},
    @{Row=3;  A=45695.05555555555; B=24.9},
    @{Row=4;  A=45696.50763888889; B=27.13},
    @{Row=5;  A=45697.61527777778; B=25.3},
    @{Row=6;  A=45700.90347222222; B=24.7},
    @{Row=7;  A=45701.41736111111; B=29.2},
    @{Row=8;  A=45705.8375;        B=25.43},
    @{Row=9;  A=45709.98125;       B=26.52},
    @{Row=10; A=45711.92152777778; B=24.77},
    @{Row=11; A=45712.14236111111; B=29.92}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
